$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All Price (D) cells are plain text in the source data (e.g. "29.936.38" uses
# "." as a thousands separator). Force text format before assigning so Excel
# does not reinterpret the strings as numbers/dates.

# --- Reorder rows 21-23 (Chainlink / Dai / WrappedliquidstakedEther2.0) ---
$ws.Range("B21").Value = 'Chainlink'
$ws.Range("C21").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.165'
$ws.Range("E21").Value = '  -1.27%  '

$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9996'
$ws.Range("E22").Value = '  -0.22%  '

$ws.Range("B23").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C23").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.127.78'
$ws.Range("E23").Value = '  -2.40%  '

# --- Price / Volume updates ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.936.38'
$ws.Range("E2").Value = '  +0.13%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.895.84'
$ws.Range("E3").Value = '  -0.06%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7751'
$ws.Range("E5").Value = '  -2.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '244.81'
$ws.Range("E6").Value = '  +0.43%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  -0.56%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '25.80'
$ws.Range("E9").Value = '  +1.83%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07256'
$ws.Range("E10").Value = '  +1.44%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08857'
$ws.Range("E11").Value = '  +9.32%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7754'
$ws.Range("E12").Value = '  +1.08%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.464'
$ws.Range("E13").Value = '  -2.07%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '94.88'
$ws.Range("E14").Value = '  +2.51%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.847.06'
$ws.Range("E15").Value = '  -3.00%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.201'
$ws.Range("E16").Value = '  +0.67%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '29.842.44'
$ws.Range("E17").Value = '  -0.29%  '
$ws.Range("E18").Value = '  +0.52%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '246.98'
$ws.Range("E19").Value = '  +1.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007916'
$ws.Range("E20").Value = '  +1.80%  '
$ws.Range("E25").Value = '  -4.91%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.566'
$ws.Range("E26").Value = '  +1.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '163.11'
$ws.Range("E27").Value = '  -0.54%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.86'
$ws.Range("E28").Value = '  +0.75%  '
$ws.Range("E29").Value = '  -0.45%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.428'
$ws.Range("E30").Value = '  +1.90%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.551'
$ws.Range("E31").Value = '  +0.04%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.536'
$ws.Range("E32").Value = '  +1.09%  '
$ws.Range("E33").Value = '  +0.81%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05528'
$ws.Range("E34").Value = '  -1.28%  '
$ws.Range("E35").Value = '  -2.29%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7562'
$ws.Range("E36").Value = '  +1.94%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9997'
$ws.Range("E37").Value = '  -0.21%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.725'
$ws.Range("E38").Value = '  +3.46%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01970'
$ws.Range("E39").Value = '  +1.92%  '
$ws.Range("E40").Value = '  +0.29%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4530'
$ws.Range("E41").Value = '  +2.45%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '74.12'
$ws.Range("E42").Value = '  -0.35%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.084'
$ws.Range("E43").Value = '  +2.27%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.090.57'
$ws.Range("E44").Value = '  -5.95%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8566'
$ws.Range("E45").Value = '  +0.50%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.000'
$ws.Range("E46").Value = '  -0.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.897'
$ws.Range("E47").Value = '  +0.84%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '102.88'
$ws.Range("E48").Value = '  -1.86%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.625'
$ws.Range("E49").Value = '  +2.29%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.887'
$ws.Range("E50").Value = '  -0.77%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.000'
$ws.Range("E51").Value = '  +0.04%  '
